$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Flip C5:C7 (turn diameter) to negative values for right turns.
$ws.Range("C5").Value = -3.82
$ws.Range("C6").Value = -2.5
$ws.Range("C7").Value = -1.92

# Move the selection to X29 (matches the saved view state in the target file).
$ws.Range("X29").Select()
